$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10:D10").NumberFormat = "@"
$ws.Range("B10").Value = "73.62"
$ws.Range("C10").Value = "23.74"
$ws.Range("D10").Value = "97.36"
